# 2018-11-12 update: the "含税/未税" (tax-inclusive/exclusive) column is being
# dropped from the cost-inquiry template. It lived in column H with its own
# border + font styling; removing the whole column shifts the trailing
# "交货日期" (delivery date) column from I back into H.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(8).Delete() | Out-Null

# Leave the selection where the editor last left it.
$ws.Range("F11").Select() | Out-Null
